# Update row 21 (2025Q3) metrics on the active sheet to reflect
# refreshed totals: total_customers=12, returning_customers=11,
# new_customers=1, recurrence_rate=3.160919540229886

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 11
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 3.160919540229886
